$wb = $excel.ActiveWorkbook

# Sheet1: add a second column of test data (invalidProgramId / 209876),
# mirroring the existing invalidBatchId / 345678 pair in column A.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Value = "invalidProgramId"
$ws1.Range("B2").Value = "'209876"
$ws1.Columns.Item(2).ColumnWidth = 120 / 7

# Make Sheet1 the active sheet/selection (previously Sheet4 was active).
$ws1.Activate()
$ws1.Range("C5").Select()
